$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 6697.6665
$ws.Range("J16").Value = 9996.5
$ws.Range("L16").Value = 9996.5
$ws.Range("N16").Value = -10456.5

$ws.Range("H28").Value = 1708
$ws.Range("I28").Value = 1664.6
$ws.Range("K28").Value = 1664.6
$ws.Range("M28").Value = -1179.6

$ws.Range("H92").Value = 980.13043
$ws.Range("I92").Value = 836.0625
$ws.Range("K92").Value = 836.0625
$ws.Range("M92").Value = 411.9375

$ws.Range("H100").Value = 5767.364
$ws.Range("I100").Value = 3376.077
$ws.Range("K100").Value = 3376.077
$ws.Range("M100").Value = -2835.077

$ws.Range("H103").Value = 1699
$ws.Range("I103").Value = 399
$ws.Range("K103").Value = 1197
$ws.Range("M103").Value = -611

$ws.Range("H138").Value = 2230.8472
$ws.Range("I138").Value = 1230.2727
$ws.Range("J138").Value = 2671.1
$ws.Range("K138").Value = 3690.8181
$ws.Range("L138").Value = 8013.299999999999
$ws.Range("M138").Value = 1449.1819
$ws.Range("N138").Value = -18293.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8088.2056
$ws.Range("I32").Value = 6421.4126
$ws.Range("K32").Value = 6421.4126
$ws.Range("M32").Value = -6134.4126

$ws.Range("H74").Value = 2686.7
$ws.Range("I74").Value = 2742.353
$ws.Range("J74").Value = 2371.3333
$ws.Range("K74").Value = 2742.353
$ws.Range("L74").Value = 2371.3333
$ws.Range("M74").Value = -1868.353
$ws.Range("N74").Value = -4119.3333

$ws.Range("H77").Value = 2686.7
$ws.Range("I77").Value = 2742.353
$ws.Range("J77").Value = 2371.3333
$ws.Range("K77").Value = 13711.765
$ws.Range("L77").Value = 11856.6665
$ws.Range("M77").Value = -9343.764999999999
$ws.Range("N77").Value = -20592.6665

$ws.Range("H97").Value = 1875.7693
$ws.Range("I97").Value = 1865.1111
$ws.Range("J97").Value = 1899.75
$ws.Range("K97").Value = 1865.1111
$ws.Range("L97").Value = 1899.75
$ws.Range("M97").Value = -1369.1111
$ws.Range("N97").Value = -2891.75

$ws.Range("H132").Value = 4225.15
$ws.Range("I132").Value = 3792.5386
$ws.Range("K132").Value = 11377.6158
$ws.Range("M132").Value = -8847.6158

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 33584.617
$ws.Range("I99").Value = 42498
$ws.Range("K99").Value = 42498
$ws.Range("M99").Value = -41000

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 18
$ws.Range("I14").Value = 18
$ws.Range("K14").Value = 18
$ws.Range("M14").Value = 152

$ws.Range("H21").Value = 2657.5
$ws.Range("J21").Value = 2657.5
$ws.Range("L21").Value = 2657.5
$ws.Range("N21").Value = -3127.5

$ws.Range("H31").Value = 3402.158
$ws.Range("I31").Value = 1322.8889
$ws.Range("K31").Value = 1322.8889
$ws.Range("M31").Value = -1027.8889

$ws.Range("H34").Value = 3402.158
$ws.Range("I34").Value = 1322.8889
$ws.Range("K34").Value = 1322.8889
$ws.Range("M34").Value = -1120.8889

$ws.Range("H58").Value = 3169.875
$ws.Range("I58").Value = 2670.889
$ws.Range("J58").Value = 3811.4285
$ws.Range("K58").Value = 2670.889
$ws.Range("L58").Value = 3811.4285
$ws.Range("M58").Value = -2467.889
$ws.Range("N58").Value = -4217.4285

$ws.Range("H86").Value = 36863.445
$ws.Range("I86").Value = 53464.332
$ws.Range("K86").Value = 53464.332
$ws.Range("M86").Value = -52341.332

$ws.Range("H89").Value = 36863.445
$ws.Range("I89").Value = 53464.332
$ws.Range("K89").Value = 267321.66
$ws.Range("M89").Value = -261705.66

$ws.Range("H99").Value = 13207175
$ws.Range("I99").Value = 4072224
$ws.Range("K99").Value = 4072224
$ws.Range("M99").Value = -4070726

$ws.Range("H126").Value = 13207175
$ws.Range("I126").Value = 4072224
$ws.Range("K126").Value = 12216672
$ws.Range("M126").Value = -12214202

$ws.Range("H132").Value = 2572.0645
$ws.Range("I132").Value = 2562.1785
$ws.Range("K132").Value = 7686.5355
$ws.Range("M132").Value = -5156.5355

$ws.Range("H134").Value = 4580.2573
$ws.Range("I134").Value = 3290.25
$ws.Range("J134").Value = 6300.2666
$ws.Range("K134").Value = 9870.75
$ws.Range("L134").Value = 18900.7998
$ws.Range("M134").Value = -7335.75
$ws.Range("N134").Value = -23970.7998

$ws.Range("H136").Value = 3169.875
$ws.Range("I136").Value = 2670.889
$ws.Range("J136").Value = 3811.4285
$ws.Range("K136").Value = 8012.667
$ws.Range("L136").Value = 11434.2855
$ws.Range("M136").Value = -5462.667
$ws.Range("N136").Value = -16534.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 826.13043
$ws.Range("I7").Value = 972.2222
$ws.Range("K7").Value = 2916.6666
$ws.Range("M7").Value = -2804.6666

$ws.Range("H63").Value = 5839.4614
$ws.Range("I63").Value = 4732
$ws.Range("J63").Value = 6040.8184
$ws.Range("K63").Value = 14196
$ws.Range("L63").Value = 18122.4552
$ws.Range("M63").Value = -13447
$ws.Range("N63").Value = -19620.4552

$ws.Range("H66").Value = 5839.4614
$ws.Range("I66").Value = 4732
$ws.Range("J66").Value = 6040.8184
$ws.Range("K66").Value = 42588
$ws.Range("L66").Value = 54367.3656
$ws.Range("M66").Value = -38844
$ws.Range("N66").Value = -61855.3656

$ws.Range("H80").Value = 4389.4
$ws.Range("I80").Value = 1550
$ws.Range("J80").Value = 4704.8887
$ws.Range("K80").Value = 4650
$ws.Range("L80").Value = 14114.6661
$ws.Range("M80").Value = -3714
$ws.Range("N80").Value = -15986.6661

$ws.Range("H83").Value = 4389.4
$ws.Range("I83").Value = 1550
$ws.Range("J83").Value = 4704.8887
$ws.Range("K83").Value = 13950
$ws.Range("L83").Value = 42343.99830000001
$ws.Range("M83").Value = -9270
$ws.Range("N83").Value = -51703.99830000001

$ws.Range("H117").Value = 4314.2856
$ws.Range("I117").Value = 700
$ws.Range("J117").Value = 4916.6665
$ws.Range("K117").Value = 2100
$ws.Range("L117").Value = 14749.9995
$ws.Range("M117").Value = 1342
$ws.Range("N117").Value = -21633.9995

$ws.Range("H131").Value = 6056.778
$ws.Range("J131").Value = 8353.4
$ws.Range("L131").Value = 25060.2
$ws.Range("N131").Value = -35140.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 2001032.8
$ws.Range("I3").Value = 3333857.8
$ws.Range("J3").Value = 1429822.1
$ws.Range("K3").Value = 3333857.8
$ws.Range("L3").Value = 1429822.1
$ws.Range("M3").Value = -3333741.8
$ws.Range("N3").Value = -1430054.1

$ws.Range("H7").Value = 4357400
$ws.Range("J7").Value = 4357400
$ws.Range("L7").Value = 4357400
$ws.Range("N7").Value = -4357624

$ws.Range("H8").Value = 4357400
$ws.Range("J8").Value = 4357400
$ws.Range("L8").Value = 4357400
$ws.Range("N8").Value = -4357678

$ws.Range("H11").Value = 1849727.2
$ws.Range("I11").Value = 2135889
$ws.Range("J11").Value = 562000
$ws.Range("K11").Value = 2135889
$ws.Range("L11").Value = 562000
$ws.Range("M11").Value = -2135750
$ws.Range("N11").Value = -562278

$ws.Range("H13").Value = 330.9
$ws.Range("J13").Value = 549.8
$ws.Range("L13").Value = 549.8
$ws.Range("N13").Value = -827.8

$ws.Range("H97").Value = 906.6667
$ws.Range("I97").Value = 1010
$ws.Range("J97").Value = 700
$ws.Range("K97").Value = 1010
$ws.Range("L97").Value = 700
$ws.Range("M97").Value = -514
$ws.Range("N97").Value = -1692

$ws.Range("H107").Value = 656
$ws.Range("I107").Value = 349.5
$ws.Range("K107").Value = 349.5
$ws.Range("M107").Value = 1570.5

$ws.Range("H126").Value = 2805.8462
$ws.Range("I126").Value = 1937.6
$ws.Range("K126").Value = 5812.799999999999
$ws.Range("M126").Value = -3342.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2502014.5
$ws.Range("J93").Value = 9999998
$ws.Range("L93").Value = 9999998
$ws.Range("N93").Value = -10002494

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 208199.2
$ws.Range("J4").Value = 10249.25
$ws.Range("L4").Value = 10249.25
$ws.Range("N4").Value = -10475.25

$ws.Range("H81").Value = 10155.917
$ws.Range("I81").Value = 21714.2
$ws.Range("J81").Value = 1900
$ws.Range("K81").Value = 43428.4
$ws.Range("L81").Value = 3800
$ws.Range("M81").Value = -42367.4
$ws.Range("N81").Value = -5922

$ws.Range("H84").Value = 10155.917
$ws.Range("I84").Value = 21714.2
$ws.Range("J84").Value = 1900
$ws.Range("K84").Value = 217142
$ws.Range("L84").Value = 19000
$ws.Range("M84").Value = -211838
$ws.Range("N84").Value = -29608

$ws.Range("H96").Value = 3783.1667
$ws.Range("I96").Value = 3500
$ws.Range("J96").Value = 3924.75
$ws.Range("K96").Value = 3500
$ws.Range("L96").Value = 3924.75
$ws.Range("M96").Value = -2127
$ws.Range("N96").Value = -6670.75

$ws.Range("H113").Value = 382.55554
$ws.Range("I113").Value = 322.73334
$ws.Range("K113").Value = 968.20002
$ws.Range("M113").Value = 1201.79998

Write-Host "Applied all cell updates"